$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "58.625.47"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.160.80"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.442"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("E11").Value = "  +4.58%  "
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("D13").Value = "3.705.23"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.80%  "
$ws.Range("D16").Value = "58.689.27"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("E17").Value = "  +4.50%  "
$ws.Range("D18").Value = "3.154.96"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "377.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.66%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  +5.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0694"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").Value = "2.651.30"
$ws.Range("E41").Value = "  +6.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.69%  "
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("E45").Value = "  +7.62%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +15.92%  "
$ws.Range("D48").Value = "3.201.45"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.978"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.67%  "
